$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.144.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.900.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'307.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.5232"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "'0.3807"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07291"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "'21.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'0.9025"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'0.08173"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "'95.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'1.848.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "'0.000008655"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'27.181.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'5.123"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("D22").Value = "'10.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "'6.456"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("B24").Value = "LidoDAOToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D24").Value = "'2.326"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'149.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "'18.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'1.741"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "'115.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "
$ws.Range("D29").Value = "'4.822"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'4.898"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'0.09224"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'0.7927"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "'1.220"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "'2.969"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'3.362"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'2.639"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "'0.5713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.01994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "'1.081"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "'9.068"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "'6.607"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "'116.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "'0.1516"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'0.4892"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'1.636"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").Value = "'38.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'0.05957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
